# Module6Output_20240103.xlsx edit script
# - DeliveryPlan (sheet1): populate header row + two data rows
# - VehicleLog (sheet2): append one data row (headers unchanged)
# - TruckUsageLog (sheet3): append one data row (headers unchanged)
# - ValidationLog (sheet5) and the remaining sheets are untouched (their
#   apparent shared-string index shifts in the diff come purely from the
#   new strings introduced on the sheets above; Excel manages that table
#   automatically as soon as we write the new cell values).

$wb = $excel.ActiveWorkbook

$wsDeliveryPlan   = $wb.Worksheets.Item("DeliveryPlan")
$wsVehicleLog     = $wb.Worksheets.Item("VehicleLog")
$wsTruckUsageLog  = $wb.Worksheets.Item("TruckUsageLog")

# ---------------------------------------------------------------------
# DeliveryPlan (sheet1)
# ---------------------------------------------------------------------
$deliveryHeaders = @(
    "vehicle_uid",
    "ori_deployment_uid",
    "material",
    "sending",
    "receiving",
    "planned_deployment_date",
    "actual_ship_date",
    "actual_delivery_date",
    "delivery_qty",
    "truck_type",
    "truck_load_pct",
    "WFR",
    "VFR"
)
for ($i = 0; $i -lt $deliveryHeaders.Length; $i++) {
    $wsDeliveryPlan.Cells.Item(1, $i + 1).Value = $deliveryHeaders[$i]
}
# Reuse the existing bold/border header style (already present on the
# other sheets) instead of building a brand-new style from scratch.
$wsVehicleLog.Range("A1").Copy()
$wsDeliveryPlan.Range("A1:M1").PasteSpecial(-4122)

$deliveryRow2 = @(
    "20240103-PLANT_001-DC_001-LARGE-#1",
    "MAT_B|PLANT_001|DC_001|2024-01-06|net demand for forecast|000046",
    "MAT_B",
    "PLANT_001",
    "DC_001"
)
for ($i = 0; $i -lt $deliveryRow2.Length; $i++) {
    $wsDeliveryPlan.Cells.Item(2, $i + 1).Value = $deliveryRow2[$i]
}
$wsDeliveryPlan.Cells.Item(2, 6).Value = 45297
$wsDeliveryPlan.Cells.Item(2, 7).Value = 45294
$wsDeliveryPlan.Cells.Item(2, 8).Value = 45296
$wsDeliveryPlan.Cells.Item(2, 9).Value = 35
$wsDeliveryPlan.Cells.Item(2, 10).Value = "LARGE"
$wsDeliveryPlan.Cells.Item(2, 11).Value = 0.76
$wsDeliveryPlan.Cells.Item(2, 12).Value = 0.76
$wsDeliveryPlan.Cells.Item(2, 13).Value = 0.7125

$deliveryRow3 = @(
    "20240103-PLANT_001-DC_001-LARGE-#1",
    "MAT_B|PLANT_001|DC_001|2024-01-01|net demand for safety|000023",
    "MAT_B",
    "PLANT_001",
    "DC_001"
)
for ($i = 0; $i -lt $deliveryRow3.Length; $i++) {
    $wsDeliveryPlan.Cells.Item(3, $i + 1).Value = $deliveryRow3[$i]
}
$wsDeliveryPlan.Cells.Item(3, 6).Value = 45292
$wsDeliveryPlan.Cells.Item(3, 7).Value = 45294
$wsDeliveryPlan.Cells.Item(3, 8).Value = 45296
$wsDeliveryPlan.Cells.Item(3, 9).Value = 60
$wsDeliveryPlan.Cells.Item(3, 10).Value = "LARGE"
$wsDeliveryPlan.Cells.Item(3, 11).Value = 0.76
$wsDeliveryPlan.Cells.Item(3, 12).Value = 0.76
$wsDeliveryPlan.Cells.Item(3, 13).Value = 0.7125

# Apply the custom date-time number format to the three date columns
# (planned_deployment_date, actual_ship_date, actual_delivery_date).
$wsDeliveryPlan.Range("F2:H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# VehicleLog (sheet2) - headers already correct, add the data row
# ---------------------------------------------------------------------
$wsVehicleLog.Cells.Item(2, 1).Value = 45294
$wsVehicleLog.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsVehicleLog.Cells.Item(2, 2).Value = "PLANT_001"
$wsVehicleLog.Cells.Item(2, 3).Value = "DC_001"
$wsVehicleLog.Cells.Item(2, 4).Value = "LARGE"
$wsVehicleLog.Cells.Item(2, 5).Value = 1
$wsVehicleLog.Cells.Item(2, 6).Value = "20240103-PLANT_001-DC_001-LARGE-#1"
$wsVehicleLog.Cells.Item(2, 7).Value = 95
$wsVehicleLog.Cells.Item(2, 8).Value = 76
$wsVehicleLog.Cells.Item(2, 9).Value = 142.5
$wsVehicleLog.Cells.Item(2, 10).Value = 0.76
$wsVehicleLog.Cells.Item(2, 11).Value = 0.7125
$wsVehicleLog.Cells.Item(2, 12).Value = "threshold"

# ---------------------------------------------------------------------
# TruckUsageLog (sheet3) - headers already correct, add the data row
# ---------------------------------------------------------------------
$wsTruckUsageLog.Cells.Item(2, 1).Value = 45294
$wsTruckUsageLog.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsTruckUsageLog.Cells.Item(2, 2).Value = "PLANT_001"
$wsTruckUsageLog.Cells.Item(2, 3).Value = "DC_001"
$wsTruckUsageLog.Cells.Item(2, 4).Value = "LARGE"
$wsTruckUsageLog.Cells.Item(2, 5).Value = 1
